$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ntrk3"
$ws.Cells.Item(2, 3).Value = "Ptprf"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.009008
$ws.Cells.Item(2, 8).Value = 0.027024
$ws.Cells.Item(2, 9).Value = 0.009775433435787767
$ws.Cells.Item(2, 10).Value = 0.009775433435787767
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.4440536666666666
$ws.Cells.Item(2, 14).Value = 1.332161
$ws.Cells.Item(2, 15).Value = 0.06598745121877762
$ws.Cells.Item(2, 16).Value = 0.06598745121877764
$ws.Cells.Item(2, 17).Value = 0.004000035429333334
$ws.Cells.Item(2, 18).Value = 0.036000318864
$ws.Cells.Item(2, 19).Value = 0.000645055936986453
$ws.Cells.Item(2, 20).Value = 0.0006450559369864531

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ntrk3"
$ws.Cells.Item(3, 3).Value = "Ptprf"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.009008
$ws.Cells.Item(3, 8).Value = 0.027024
$ws.Cells.Item(3, 9).Value = 0.009775433435787767
$ws.Cells.Item(3, 10).Value = 0.009775433435787767
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.674351333333334
$ws.Cells.Item(3, 14).Value = 11.023054
$ws.Cells.Item(3, 15).Value = 0.5460175144798202
$ws.Cells.Item(3, 16).Value = 0.5460175144798202
$ws.Cells.Item(3, 17).Value = 0.03309855681066667
$ws.Cells.Item(3, 18).Value = 0.297887011296
$ws.Cells.Item(3, 19).Value = 0.005337557867571765
$ws.Cells.Item(3, 20).Value = 0.005337557867571765

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ntrk3"
$ws.Cells.Item(4, 3).Value = "Ptprf"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.009008
$ws.Cells.Item(4, 8).Value = 0.027024
$ws.Cells.Item(4, 9).Value = 0.009775433435787767
$ws.Cells.Item(4, 10).Value = 0.009775433435787767
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.008309666666666667
$ws.Cells.Item(4, 14).Value = 0.024929
$ws.Cells.Item(4, 15).Value = 0.001234836608662848
$ws.Cells.Item(4, 16).Value = 0.001234836608662848
$ws.Cells.Item(4, 17).Value = 0.00007485347733333334
$ws.Cells.Item(4, 18).Value = 0.000673681296
$ws.Cells.Item(4, 19).Value = 0.00001207106307205757
$ws.Cells.Item(4, 20).Value = 0.00001207106307205757

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Ntrk3"
$ws.Cells.Item(5, 3).Value = "Ptprf"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.009008
$ws.Cells.Item(5, 8).Value = 0.027024
$ws.Cells.Item(5, 9).Value = 0.009775433435787767
$ws.Cells.Item(5, 10).Value = 0.009775433435787767
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.602650666666667
$ws.Cells.Item(5, 14).Value = 7.807952
$ws.Cells.Item(5, 15).Value = 0.3867601976927393
$ws.Cells.Item(5, 16).Value = 0.3867601976927393
$ws.Cells.Item(5, 17).Value = 0.02344467720533334
$ws.Cells.Item(5, 18).Value = 0.211002094848
$ws.Cells.Item(5, 19).Value = 0.003780748568157491
$ws.Cells.Item(5, 20).Value = 0.00378074856815749

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ntrk3"
$ws.Cells.Item(6, 3).Value = "Ptprf"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.07403633333333333
$ws.Cells.Item(6, 8).Value = 0.222109
$ws.Cells.Item(6, 9).Value = 0.08034383307391152
$ws.Cells.Item(6, 10).Value = 0.08034383307391152
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.4440536666666666
$ws.Cells.Item(6, 14).Value = 1.332161
$ws.Cells.Item(6, 15).Value = 0.06598745121877762
$ws.Cells.Item(6, 16).Value = 0.06598745121877764
$ws.Cells.Item(6, 17).Value = 0.03287610528322222
$ws.Cells.Item(6, 18).Value = 0.295884947549
$ws.Cells.Item(6, 19).Value = 0.005301684765694348
$ws.Cells.Item(6, 20).Value = 0.00530168476569435

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ntrk3"
$ws.Cells.Item(7, 3).Value = "Ptprf"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.07403633333333333
$ws.Cells.Item(7, 8).Value = 0.222109
$ws.Cells.Item(7, 9).Value = 0.08034383307391152
$ws.Cells.Item(7, 10).Value = 0.08034383307391152
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.674351333333334
$ws.Cells.Item(7, 14).Value = 11.023054
$ws.Cells.Item(7, 15).Value = 0.5460175144798202
$ws.Cells.Item(7, 16).Value = 0.5460175144798202
$ws.Cells.Item(7, 17).Value = 0.2720355000984445
$ws.Cells.Item(7, 18).Value = 2.448319500886
$ws.Cells.Item(7, 19).Value = 0.04386914003879874
$ws.Cells.Item(7, 20).Value = 0.04386914003879874

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Ntrk3"
$ws.Cells.Item(8, 3).Value = "Ptprf"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.07403633333333333
$ws.Cells.Item(8, 8).Value = 0.222109
$ws.Cells.Item(8, 9).Value = 0.08034383307391152
$ws.Cells.Item(8, 10).Value = 0.08034383307391152
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.008309666666666667
$ws.Cells.Item(8, 14).Value = 0.024929
$ws.Cells.Item(8, 15).Value = 0.001234836608662848
$ws.Cells.Item(8, 16).Value = 0.001234836608662848
$ws.Cells.Item(8, 17).Value = 0.0006152172512222222
$ws.Cells.Item(8, 18).Value = 0.005536955261
$ws.Cells.Item(8, 19).Value = 0.00009921150635996283
$ws.Cells.Item(8, 20).Value = 0.00009921150635996283

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Ntrk3"
$ws.Cells.Item(9, 3).Value = "Ptprf"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.07403633333333333
$ws.Cells.Item(9, 8).Value = 0.222109
$ws.Cells.Item(9, 9).Value = 0.08034383307391152
$ws.Cells.Item(9, 10).Value = 0.08034383307391152
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.602650666666667
$ws.Cells.Item(9, 14).Value = 7.807952
$ws.Cells.Item(9, 15).Value = 0.3867601976927393
$ws.Cells.Item(9, 16).Value = 0.3867601976927393
$ws.Cells.Item(9, 17).Value = 0.1926907123075556
$ws.Cells.Item(9, 18).Value = 1.734216410768
$ws.Cells.Item(9, 19).Value = 0.03107379676305847
$ws.Cells.Item(9, 20).Value = 0.03107379676305846

$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Ntrk3"
$ws.Cells.Item(10, 3).Value = "Ptprf"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.8384493333333333
$ws.Cells.Item(10, 8).Value = 2.515348
$ws.Cells.Item(10, 9).Value = 0.9098807334903006
$ws.Cells.Item(10, 10).Value = 0.9098807334903007
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.4440536666666666
$ws.Cells.Item(10, 14).Value = 1.332161
$ws.Cells.Item(10, 15).Value = 0.06598745121877762
$ws.Cells.Item(10, 16).Value = 0.06598745121877764
$ws.Cells.Item(10, 17).Value = 0.3723165007808888
$ws.Cells.Item(10, 18).Value = 3.350848507028
$ws.Cells.Item(10, 19).Value = 0.06004071051609682
$ws.Cells.Item(10, 20).Value = 0.06004071051609684

$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Ntrk3"
$ws.Cells.Item(11, 3).Value = "Ptprf"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.8384493333333333
$ws.Cells.Item(11, 8).Value = 2.515348
$ws.Cells.Item(11, 9).Value = 0.9098807334903006
$ws.Cells.Item(11, 10).Value = 0.9098807334903007
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 3.674351333333334
$ws.Cells.Item(11, 14).Value = 11.023054
$ws.Cells.Item(11, 15).Value = 0.5460175144798202
$ws.Cells.Item(11, 16).Value = 0.5460175144798202
$ws.Cells.Item(11, 17).Value = 3.080757425865777
$ws.Cells.Item(11, 18).Value = 27.726816832792
$ws.Cells.Item(11, 19).Value = 0.4968108165734496
$ws.Cells.Item(11, 20).Value = 0.4968108165734497

$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Ntrk3"
$ws.Cells.Item(12, 3).Value = "Ptprf"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.8384493333333333
$ws.Cells.Item(12, 8).Value = 2.515348
$ws.Cells.Item(12, 9).Value = 0.9098807334903006
$ws.Cells.Item(12, 10).Value = 0.9098807334903007
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.008309666666666667
$ws.Cells.Item(12, 14).Value = 0.024929
$ws.Cells.Item(12, 15).Value = 0.001234836608662848
$ws.Cells.Item(12, 16).Value = 0.001234836608662848
$ws.Cells.Item(12, 17).Value = 0.006967234476888888
$ws.Cells.Item(12, 18).Value = 0.062705110292
$ws.Cells.Item(12, 19).Value = 0.001123554039230827
$ws.Cells.Item(12, 20).Value = 0.001123554039230827

$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Ntrk3"
$ws.Cells.Item(13, 3).Value = "Ptprf"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.8384493333333333
$ws.Cells.Item(13, 8).Value = 2.515348
$ws.Cells.Item(13, 9).Value = 0.9098807334903006
$ws.Cells.Item(13, 10).Value = 0.9098807334903007
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 2.602650666666667
$ws.Cells.Item(13, 14).Value = 7.807952
$ws.Cells.Item(13, 15).Value = 0.3867601976927393
$ws.Cells.Item(13, 16).Value = 0.3867601976927393
$ws.Cells.Item(13, 17).Value = 2.182190716366222
$ws.Cells.Item(13, 18).Value = 19.639716447296
$ws.Cells.Item(13, 19).Value = 0.3519056523615233
$ws.Cells.Item(13, 20).Value = 0.3519056523615233
